$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.180180072784424
$ws.Range("B1").Value = 2.077945470809937
$ws.Range("C1").Value = 5.655363082885742
$ws.Range("D1").Value = 0.8498803973197937
$ws.Range("E1").Value = 1.022581934928894
